$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# stmt sheet: insert new node "GotoStmt" (alphabetically between "ForStmt"
# and "IfStmt"), marked as completed ("o"). Also flip the existing
# "LabelStmt" and "WrapperStmt" rows to completed ("o") - they previously had
# no status.
# (Added before "LabelDecl" below so the shared-string table picks up
# "GotoStmt" first, matching the target workbook's string ordering.)
# ---------------------------------------------------------------------------
$stmtWs = $wb.Worksheets.Item("stmt")
$stmtWs.Rows.Item(17).Insert()
$stmtWs.Range("A17").Value = "GotoStmt"
$stmtWs.Range("B17").Value = "o"

# ---------------------------------------------------------------------------
# -decl- sheet: insert new node "LabelDecl" (alphabetically between
# "IndirectFieldDecl" and "LinkageSpecDecl"), marked as completed ("o").
# ---------------------------------------------------------------------------
$declWs = $wb.Worksheets.Item("-decl-")
$declWs.Rows.Item(21).Insert()
$declWs.Range("A21").Value = "LabelDecl"
$declWs.Range("B21").Value = "o"

$stmtWs.Range("B19").Value = "o"   # LabelStmt (now on row 19 after insert)
$stmtWs.Range("B28").Value = "o"   # WrapperStmt (now on row 28 after insert)

# ---------------------------------------------------------------------------
# Update the active selection on the affected sheets to match the target
# workbook. "-decl-" is selected (and thus briefly made active) first, then
# "stmt" is selected last so that it ends up as the active sheet/tab again,
# same as in the source workbook (Summary/-decl-/-expr-/stmt/-type-/-attr-
# /other, with "stmt" originally active).
# ---------------------------------------------------------------------------
$declWs.Range("B22").Select()
$stmtWs.Range("B3").Select()
